$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 154738
$ws.Range("C4").Value = 145868
$ws.Range("C7").Value = 5.73
$ws.Range("C8").Value = 63.56
